# feat: add 2022-Q3 data
#
# - Insert a new worksheet "2022-Q3" right before the existing "2022-Q2"
#   sheet (so the tab order becomes 总计, 2022-Q3, 2022-Q2) and populate it
#   with the Q3 fund-holding rows.
# - On the "总计" (totals) summary sheet, push the existing 2022-Q2 totals
#   row down to row 3 and write the new 2022-Q3 totals into row 2.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item("总计")
$q2 = $wb.Worksheets.Item("2022-Q2")

# ---------------------------------------------------------------------
# 1. New "2022-Q3" sheet, inserted before the "2022-Q2" sheet.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($q2)
$q3.Name = "2022-Q3"

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Header row + the row-index column (A) use the same bold / centered /
# thin-bordered look as the other sheets in this workbook.
$q3.Range("B1:H1").Font.Bold = $true
$q3.Range("B1:H1").HorizontalAlignment = -4108
$q3.Range("B1:H1").VerticalAlignment = -4160
$q3.Range("B1:H1").Borders.LineStyle = 1

$q3.Range("A2:A4").Font.Bold = $true
$q3.Range("A2:A4").HorizontalAlignment = -4108
$q3.Range("A2:A4").VerticalAlignment = -4160
$q3.Range("A2:A4").Borders.LineStyle = 1

# Columns B:G hold plain text in the source data (fund code / name / size
# / position figures), not numbers - e.g. the fund code "012098" must
# keep its leading zero, and "12.65" etc. must stay text, not become a
# number. A leading apostrophe is the standard Excel way to force text
# entry without leaving a lasting NumberFormat override on the cell.
$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "'012098"
$q3.Range("C2").Value = "华夏成长机会一年持有期混合"
$q3.Range("D2").Value = "'12.65"
$q3.Range("E2").Value = "'71.62"
$q3.Range("F2").Value = "'4.06"
$q3.Range("G2").Value = "'0.5136"
$q3.Range("H2").Value = 4

$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "'012924"
$q3.Range("C3").Value = "华夏新时代灵活配置混合（QDII）美元现汇"
$q3.Range("D3").Value = "'2.11"
$q3.Range("E3").Value = "'73.45"
$q3.Range("F3").Value = "'2.74"
$q3.Range("G3").Value = "'0.0578"
$q3.Range("H3").Value = 8

$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "'012925"
$q3.Range("C4").Value = "华夏新时代灵活配置混合（QDII）美元现钞"
$q3.Range("D4").Value = "'2.11"
$q3.Range("E4").Value = "'73.45"
$q3.Range("F4").Value = "'2.74"
$q3.Range("G4").Value = "'0.0578"
$q3.Range("H4").Value = 8

# ---------------------------------------------------------------------
# 2. "总计" sheet: shift the old 2022-Q2 row down to row 3, then write
#    the new 2022-Q3 totals into row 2.
# ---------------------------------------------------------------------
$total.Rows.Item(2).Insert()
# Inserting a row copies the formatting of the row above into the blank
# row; the new row 2 here is a plain data row (like row 3), so drop that
# inherited formatting before filling it in.
$total.Range("A2:D2").ClearFormats()

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.17

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.63

# The row-insert left the new A2 without the header-row formatting that
# the rest of column A carries (bold, centered, thin border) - restore it
# to match A1/A3.
$total.Range("A2").Font.Bold = $true
$total.Range("A2").HorizontalAlignment = -4108
$total.Range("A2").VerticalAlignment = -4160
$total.Range("A2").Borders.LineStyle = 1
